# Orbital User to MAGELLANO analysis
# Update the apoapsis/periapsis altitude inputs on the "Mars-Mars" sheet
# and move the active selection to C14 (matching the scroll/selection
# state left behind by the author after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mars-Mars")

# Orbit 1 (Small): h_Apo (C7) and h_Per (C13) inputs
$ws.Range("C7").Value = 12000
$ws.Range("C13").Value = 12000

# Leave the selection on C14, scrolled back to the top (no frozen
# topLeftCell override), as in the saved workbook.
$ws.Activate()
$ws.Range("C14").Select()
